# Natmi following Dr Hou advice
# Rebuild the LR-pairs data table (Angptl3 -> Itgb3) with the updated
# cell/expression-count methodology. The sending clusters (FAPs, M2, sCs)
# now each connect to 4 target clusters (ECs, FAPs, M2, sCs) instead of 3,
# and every metric column is recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angptl3"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.081251
$ws.Range("H2").Value = 9.243753
$ws.Range("I2").Value = 0.6013642694204734
$ws.Range("J2").Value = 0.6013642694204734
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.481489333333333
$ws.Range("N2").Value = 7.444467999999999
$ws.Range("O2").Value = 0.2345069082418988
$ws.Range("P2").Value = 0.2345069082418987
$ws.Range("Q2").Value = 7.646091489822664
$ws.Range("R2").Value = 68.81482340840398
$ws.Range("S2").Value = 0.1410240755489434
$ws.Range("T2").Value = 0.1410240755489434

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angptl3"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.081251
$ws.Range("H3").Value = 9.243753
$ws.Range("I3").Value = 0.6013642694204734
$ws.Range("J3").Value = 0.6013642694204734
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.245227
$ws.Range("N3").Value = 21.735681
$ws.Range("O3").Value = 0.6846919551326144
$ws.Range("P3").Value = 0.6846919551326142
$ws.Range("Q3").Value = 22.324362938977
$ws.Range("R3").Value = 200.919266450793
$ws.Range("S3").Value = 0.4117492773764002
$ws.Range("T3").Value = 0.4117492773764

# Row 4: FAPs -> M2
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angptl3"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.081251
$ws.Range("H4").Value = 9.243753
$ws.Range("I4").Value = 0.6013642694204734
$ws.Range("J4").Value = 0.6013642694204734
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2001876666666667
$ws.Range("N4").Value = 0.600563
$ws.Range("O4").Value = 0.01891823194544989
$ws.Range("P4").Value = 0.01891823194544989
$ws.Range("Q4").Value = 0.6168284481043332
$ws.Range("R4").Value = 5.551456032939
$ws.Range("S4").Value = 0.01137674873260253
$ws.Range("T4").Value = 0.01137674873260253

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Angptl3"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.081251
$ws.Range("H5").Value = 9.243753
$ws.Range("I5").Value = 0.6013642694204734
$ws.Range("J5").Value = 0.6013642694204734
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6548283333333333
$ws.Range("N5").Value = 1.964485
$ws.Range("O5").Value = 0.06188290468003712
$ws.Range("P5").Value = 0.06188290468003711
$ws.Range("Q5").Value = 2.017690456911667
$ws.Range("R5").Value = 18.159214112205
$ws.Range("S5").Value = 0.03721416776252732
$ws.Range("T5").Value = 0.0372141677625273

# Row 6: M2 -> ECs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Angptl3"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6679959999999999
$ws.Range("H6").Value = 2.003988
$ws.Range("I6").Value = 0.1303720230892577
$ws.Range("J6").Value = 0.1303720230892577
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.481489333333333
$ws.Range("N6").Value = 7.444467999999999
$ws.Range("O6").Value = 0.2345069082418988
$ws.Range("P6").Value = 0.2345069082418987
$ws.Range("Q6").Value = 1.657624948709333
$ws.Range("R6").Value = 14.91862453838399
$ws.Range("S6").Value = 0.03057314005590327
$ws.Range("T6").Value = 0.03057314005590327

# Row 7: M2 -> FAPs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Angptl3"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6679959999999999
$ws.Range("H7").Value = 2.003988
$ws.Range("I7").Value = 0.1303720230892577
$ws.Range("J7").Value = 0.1303720230892577
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.245227
$ws.Range("N7").Value = 21.735681
$ws.Range("O7").Value = 0.6846919551326144
$ws.Range("P7").Value = 0.6846919551326142
$ws.Range("Q7").Value = 4.839782655092
$ws.Range("R7").Value = 43.55804389582799
$ws.Range("S7").Value = 0.08926467538357823
$ws.Range("T7").Value = 0.0892646753835782

# Row 8: M2 -> M2
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Angptl3"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6679959999999999
$ws.Range("H8").Value = 2.003988
$ws.Range("I8").Value = 0.1303720230892577
$ws.Range("J8").Value = 0.1303720230892577
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2001876666666667
$ws.Range("N8").Value = 0.600563
$ws.Range("O8").Value = 0.01891823194544989
$ws.Range("P8").Value = 0.01891823194544989
$ws.Range("Q8").Value = 0.1337245605826666
$ws.Range("R8").Value = 1.203521045244
$ws.Range("S8").Value = 0.002466408172000127
$ws.Range("T8").Value = 0.002466408172000126

# Row 9: M2 -> sCs
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Angptl3"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6679959999999999
$ws.Range("H9").Value = 2.003988
$ws.Range("I9").Value = 0.1303720230892577
$ws.Range("J9").Value = 0.1303720230892577
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6548283333333333
$ws.Range("N9").Value = 1.964485
$ws.Range("O9").Value = 0.06188290468003712
$ws.Range("P9").Value = 0.06188290468003711
$ws.Range("Q9").Value = 0.4374227073533333
$ws.Range("R9").Value = 3.936804366179999
$ws.Range("S9").Value = 0.008067799477776135
$ws.Range("T9").Value = 0.008067799477776133

# Row 10: sCs -> ECs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Angptl3"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.374521
$ws.Range("H10").Value = 4.123563
$ws.Range("I10").Value = 0.2682637074902688
$ws.Range("J10").Value = 0.2682637074902689
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.481489333333333
$ws.Range("N10").Value = 7.444467999999999
$ws.Range("O10").Value = 0.2345069082418988
$ws.Range("P10").Value = 0.2345069082418987
$ws.Range("Q10").Value = 3.410859199942665
$ws.Range("R10").Value = 30.69773279948399
$ws.Range("S10").Value = 0.06290969263705204
$ws.Range("T10").Value = 0.06290969263705204

# Row 11: sCs -> FAPs
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Angptl3"
$ws.Range("C11").Value = "Itgb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.374521
$ws.Range("H11").Value = 4.123563
$ws.Range("I11").Value = 0.2682637074902688
$ws.Range("J11").Value = 0.2682637074902689
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.245227
$ws.Range("N11").Value = 21.735681
$ws.Range("O11").Value = 0.6846919551326144
$ws.Range("P11").Value = 0.6846919551326142
$ws.Range("Q11").Value = 9.958716661266999
$ws.Range("R11").Value = 89.628449951403
$ws.Range("S11").Value = 0.183678002372636
$ws.Range("T11").Value = 0.1836780023726359

# Row 12: sCs -> M2
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Angptl3"
$ws.Range("C12").Value = "Itgb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.374521
$ws.Range("H12").Value = 4.123563
$ws.Range("I12").Value = 0.2682637074902688
$ws.Range("J12").Value = 0.2682637074902689
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2001876666666667
$ws.Range("N12").Value = 0.600563
$ws.Range("O12").Value = 0.01891823194544989
$ws.Range("P12").Value = 0.01891823194544989
$ws.Range("Q12").Value = 0.2751621517743333
$ws.Range("R12").Value = 2.476459365969
$ws.Range("S12").Value = 0.005075075040847229
$ws.Range("T12").Value = 0.00507507504084723

# Row 13: sCs -> sCs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Angptl3"
$ws.Range("C13").Value = "Itgb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.374521
$ws.Range("H13").Value = 4.123563
$ws.Range("I13").Value = 0.2682637074902688
$ws.Range("J13").Value = 0.2682637074902689
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6548283333333333
$ws.Range("N13").Value = 1.964485
$ws.Range("O13").Value = 0.06188290468003712
$ws.Range("P13").Value = 0.06188290468003711
$ws.Range("Q13").Value = 0.9000752955616667
$ws.Range("R13").Value = 8.100677660055
$ws.Range("S13").Value = 0.01660093743973367
$ws.Range("T13").Value = 0.01660093743973367
